# Auto-generated Excel COM-interop edit script
# Applies: (1) shared-string text fix (mojibake -> accented chars)
#          (2) 161 recalculated numeric cell values on rows 67-73, 97-98 of sheet Tab04

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab04")

# Fix mojibake in the Regional Economic Communities legend note (cell A103)
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Países Africanos de Língua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Común del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'

# Update recalculated numeric values
$ws.Cells.Item(67, 7).Value = [double]"0.33436389165648001"
$ws.Cells.Item(67, 11).Value = [double]"2.2598254517496499"
$ws.Cells.Item(67, 12).Value = [double]"1.6395345576412701"
$ws.Cells.Item(67, 13).Value = [double]"1.8142244229721001"
$ws.Cells.Item(67, 14).Value = [double]"1.6096143935993401"
$ws.Cells.Item(67, 15).Value = [double]"-0.71435840265119999"
$ws.Cells.Item(67, 18).Value = [double]"3.1437521620717899"
$ws.Cells.Item(67, 19).Value = [double]"2.97852557562253"
$ws.Cells.Item(67, 21).Value = [double]"2.5277052213361602"
$ws.Cells.Item(67, 27).Value = [double]"8.9248514986269997E-2"
$ws.Cells.Item(67, 28).Value = [double]"1.73813365623288"
$ws.Cells.Item(67, 29).Value = [double]"1.1309363464845601"
$ws.Cells.Item(67, 31).Value = [double]"2.2506871921263398"
$ws.Cells.Item(67, 34).Value = [double]"2.9544641904655902"
$ws.Cells.Item(67, 35).Value = [double]"2.1594353925808001"
$ws.Cells.Item(67, 36).Value = [double]"2.22212558367589"
$ws.Cells.Item(67, 38).Value = [double]"3.03673184338022"
$ws.Cells.Item(67, 40).Value = [double]"3.1536938150549898"
$ws.Cells.Item(67, 41).Value = [double]"3.2442409378534101"
$ws.Cells.Item(68, 3).Value = [double]"0.34649512151114997"
$ws.Cells.Item(68, 4).Value = [double]"2.3033366459201101"
$ws.Cells.Item(68, 7).Value = [double]"0.96222254739461999"
$ws.Cells.Item(68, 8).Value = [double]"-1.3193257541636001"
$ws.Cells.Item(68, 9).Value = [double]"2.62570638188045"
$ws.Cells.Item(68, 10).Value = [double]"1.0180564704075199"
$ws.Cells.Item(68, 17).Value = [double]"3.5161245963979901"
$ws.Cells.Item(68, 18).Value = [double]"2.9985432015000999"
$ws.Cells.Item(68, 21).Value = [double]"3.1496730802139701"
$ws.Cells.Item(68, 22).Value = [double]"1.6621242880917599"
$ws.Cells.Item(68, 23).Value = [double]"4.0270002675847003"
$ws.Cells.Item(68, 24).Value = [double]"0.64386324818027996"
$ws.Cells.Item(68, 25).Value = [double]"5.1725237203738201"
$ws.Cells.Item(68, 26).Value = [double]"0.81097104250718"
$ws.Cells.Item(68, 27).Value = [double]"0.97418834040568003"
$ws.Cells.Item(68, 34).Value = [double]"2.5071668067266399"
$ws.Cells.Item(68, 35).Value = [double]"1.49601210087337"
$ws.Cells.Item(68, 38).Value = [double]"2.2403249622732599"
$ws.Cells.Item(68, 40).Value = [double]"2.3444917838142798"
$ws.Cells.Item(68, 42).Value = [double]"2.2825996676136802"
$ws.Cells.Item(69, 11).Value = [double]"0.18241313140634999"
$ws.Cells.Item(69, 12).Value = [double]"0.21725805109108001"
$ws.Cells.Item(69, 21).Value = [double]"2.15764636141316"
$ws.Cells.Item(69, 36).Value = [double]"2.6587253610310899"
$ws.Cells.Item(69, 37).Value = [double]"3.07100868739309"
$ws.Cells.Item(69, 39).Value = [double]"3.5081336301822201"
$ws.Cells.Item(69, 40).Value = [double]"3.5123067161331001"
$ws.Cells.Item(70, 5).Value = [double]"-3.0918369681171001"
$ws.Cells.Item(70, 6).Value = [double]"-3.2599649281471001"
$ws.Cells.Item(70, 8).Value = [double]"3.0775700017133101"
$ws.Cells.Item(70, 9).Value = [double]"1.85033665355485"
$ws.Cells.Item(70, 10).Value = [double]"2.3343385720305601"
$ws.Cells.Item(70, 18).Value = [double]"5.52634368226619"
$ws.Cells.Item(70, 19).Value = [double]"3.63466829709103"
$ws.Cells.Item(70, 20).Value = [double]"5.7418701436988897"
$ws.Cells.Item(70, 23).Value = [double]"1.3270214594324901"
$ws.Cells.Item(70, 25).Value = [double]"4.3360401188021296"
$ws.Cells.Item(70, 29).Value = [double]"-4.0092975283226"
$ws.Cells.Item(70, 30).Value = [double]"-2.3305249623163999"
$ws.Cells.Item(70, 32).Value = [double]"-1.2474887117613001"
$ws.Cells.Item(70, 36).Value = [double]"1.1000358358531901"
$ws.Cells.Item(70, 39).Value = [double]"1.7741835365259699"
$ws.Cells.Item(70, 40).Value = [double]"1.80761970055907"
$ws.Cells.Item(72, 7).Value = [double]"0.73500875760263995"
$ws.Cells.Item(72, 9).Value = [double]"3.4790420302464802"
$ws.Cells.Item(72, 10).Value = [double]"0.66272550877577996"
$ws.Cells.Item(72, 12).Value = [double]"1.3510655696585201"
$ws.Cells.Item(72, 22).Value = [double]"0.12787778230891"
$ws.Cells.Item(72, 23).Value = [double]"3.8874294356793202"
$ws.Cells.Item(72, 28).Value = [double]"3.7027667711180499"
$ws.Cells.Item(72, 29).Value = [double]"1.97368302478012"
$ws.Cells.Item(72, 30).Value = [double]"2.3666679728467299"
$ws.Cells.Item(72, 32).Value = [double]"2.13682487741502"
$ws.Cells.Item(72, 35).Value = [double]"1.45725013277342"
$ws.Cells.Item(73, 3).Value = [double]"-2.7935390460876999"
$ws.Cells.Item(73, 6).Value = [double]"-1.1310408084237"
$ws.Cells.Item(73, 8).Value = [double]"1.1205784511946999"
$ws.Cells.Item(73, 19).Value = [double]"3.50804614792114"
$ws.Cells.Item(73, 31).Value = [double]"-0.48876782316779999"
$ws.Cells.Item(73, 37).Value = [double]"0.81433052275619"
$ws.Cells.Item(73, 38).Value = [double]"0.76700168445943995"
$ws.Cells.Item(73, 39).Value = [double]"0.83152842205083999"
$ws.Cells.Item(97, 3).Value = [double]"-0.44734780130469998"
$ws.Cells.Item(97, 4).Value = [double]"1.09234498091094"
$ws.Cells.Item(97, 5).Value = [double]"-2.7306013596667"
$ws.Cells.Item(97, 6).Value = [double]"-1.8947285852944999"
$ws.Cells.Item(97, 7).Value = [double]"-1.0931000303572"
$ws.Cells.Item(97, 8).Value = [double]"-1.8588219451704"
$ws.Cells.Item(97, 9).Value = [double]"2.1430639466704502"
$ws.Cells.Item(97, 10).Value = [double]"-7.7009291282600006E-2"
$ws.Cells.Item(97, 11).Value = [double]"-0.17639415180260001"
$ws.Cells.Item(97, 12).Value = [double]"-0.61845599520009997"
$ws.Cells.Item(97, 13).Value = [double]"0.64769545069107004"
$ws.Cells.Item(97, 14).Value = [double]"2.2119762709934299"
$ws.Cells.Item(97, 15).Value = [double]"3.6167713906899599"
$ws.Cells.Item(97, 16).Value = [double]"4.07710214080066"
$ws.Cells.Item(97, 17).Value = [double]"4.50494823512682"
$ws.Cells.Item(97, 18).Value = [double]"4.4534364652730201"
$ws.Cells.Item(97, 19).Value = [double]"2.89627361035325"
$ws.Cells.Item(97, 20).Value = [double]"4.2034390911231299"
$ws.Cells.Item(97, 21).Value = [double]"2.7471581228378299"
$ws.Cells.Item(97, 22).Value = [double]"1.9430923352174201"
$ws.Cells.Item(97, 23).Value = [double]"5.2188271326195697"
$ws.Cells.Item(97, 24).Value = [double]"-0.59230826364069999"
$ws.Cells.Item(97, 25).Value = [double]"9.36634439659389"
$ws.Cells.Item(97, 26).Value = [double]"1.20405769497236"
$ws.Cells.Item(97, 27).Value = [double]"1.72259886369201"
$ws.Cells.Item(97, 28).Value = [double]"0.97426410871093005"
$ws.Cells.Item(97, 29).Value = [double]"-1.2499316042924999"
$ws.Cells.Item(97, 30).Value = [double]"2.1992792199702902"
$ws.Cells.Item(97, 31).Value = [double]"1.0453988046598099"
$ws.Cells.Item(97, 32).Value = [double]"0.43136672693807998"
$ws.Cells.Item(97, 33).Value = [double]"-3.9812016210591001"
$ws.Cells.Item(97, 34).Value = [double]"2.9851328628420499"
$ws.Cells.Item(97, 35).Value = [double]"0.99353592748064001"
$ws.Cells.Item(97, 36).Value = [double]"2.34114518157158"
$ws.Cells.Item(97, 37).Value = [double]"2.2227588353682801"
$ws.Cells.Item(97, 38).Value = [double]"2.2586521129128299"
$ws.Cells.Item(97, 39).Value = [double]"2.1638288238100798"
$ws.Cells.Item(97, 40).Value = [double]"2.39168122763658"
$ws.Cells.Item(97, 41).Value = [double]"2.4204366364976999"
$ws.Cells.Item(97, 42).Value = [double]"2.29142385897749"
$ws.Cells.Item(98, 3).Value = [double]"5.8140949358154597"
$ws.Cells.Item(98, 4).Value = [double]"6.13774515427676"
$ws.Cells.Item(98, 5).Value = [double]"2.7405377937209598"
$ws.Cells.Item(98, 6).Value = [double]"-1.4915076890868"
$ws.Cells.Item(98, 7).Value = [double]"-1.6671172767578"
$ws.Cells.Item(98, 8).Value = [double]"1.38361093545656"
$ws.Cells.Item(98, 9).Value = [double]"2.6099268450711199"
$ws.Cells.Item(98, 10).Value = [double]"0.65214340939227999"
$ws.Cells.Item(98, 11).Value = [double]"0.55733704619660995"
$ws.Cells.Item(98, 12).Value = [double]"-9.5827829929999997E-2"
$ws.Cells.Item(98, 13).Value = [double]"2.9637472126619402"
$ws.Cells.Item(98, 14).Value = [double]"1.16266290687814"
$ws.Cells.Item(98, 15).Value = [double]"2.09879920696625"
$ws.Cells.Item(98, 16).Value = [double]"3.5250432045696098"
$ws.Cells.Item(98, 17).Value = [double]"5.6874228881591504"
$ws.Cells.Item(98, 18).Value = [double]"4.7661083490292002"
$ws.Cells.Item(98, 19).Value = [double]"4.9064925703148097"
$ws.Cells.Item(98, 20).Value = [double]"6.0320447226581901"
$ws.Cells.Item(98, 21).Value = [double]"2.0896684401686199"
$ws.Cells.Item(98, 22).Value = [double]"-0.44849170701719998"
$ws.Cells.Item(98, 23).Value = [double]"2.6058425262479799"
$ws.Cells.Item(98, 24).Value = [double]"2.5205795979075698"
$ws.Cells.Item(98, 25).Value = [double]"0.73300107911047996"
$ws.Cells.Item(98, 26).Value = [double]"0.85872308050076995"
$ws.Cells.Item(98, 27).Value = [double]"2.52545993475852"
$ws.Cells.Item(98, 28).Value = [double]"0.49336672775056001"
$ws.Cells.Item(98, 29).Value = [double]"2.68050515833411"
$ws.Cells.Item(98, 30).Value = [double]"1.5130708166742299"
$ws.Cells.Item(98, 31).Value = [double]"0.78610571336641"
$ws.Cells.Item(98, 32).Value = [double]"-7.6778252307E-2"
$ws.Cells.Item(98, 33).Value = [double]"-0.82221423856999998"
$ws.Cells.Item(98, 34).Value = [double]"3.2159865509805901"
$ws.Cells.Item(98, 35).Value = [double]"3.5527832215207602"
$ws.Cells.Item(98, 36).Value = [double]"1.4436863904544499"
$ws.Cells.Item(98, 37).Value = [double]"2.5209613605545802"
$ws.Cells.Item(98, 38).Value = [double]"2.96489211935938"
$ws.Cells.Item(98, 39).Value = [double]"3.2033604421527699"
$ws.Cells.Item(98, 40).Value = [double]"3.36966219944365"
$ws.Cells.Item(98, 41).Value = [double]"3.35308848182992"
$ws.Cells.Item(98, 42).Value = [double]"3.08190757032179"
